# Add data for 2022-06-12
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-04"

# Update the column header label (I1) to reflect new "through" date
$ws.Range("I1").Value = "2022 (through 06-04)"

# Update June's 2022 value (I7) and the Total row (I14)
$ws.Range("I7").Value = 10
$ws.Range("I14").Value = 674
